$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 33)
$ws.Cells.Item(33, 1).Value = 10002
$ws.Cells.Item(33, 2).Value = 110032
$ws.Cells.Item(33, 3).Value = "eng"
$ws.Cells.Item(33, 4).Value = $true
$ws.Cells.Item(33, 5).Value = "superadmin"
$ws.Cells.Item(33, 6).Value = "now()"
$ws.Cells.Item(33, 7).Value = "now()"

# Update the selection to match the diff (D26)
$ws.Range("D26").Select()
